# Generate Report for handback
# - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#   (this shared string is reused by the Overview sheet too, so replace it
#   everywhere it appears so every referencing cell updates together)
# - Each per-language sheet (zh-cn / de-de) gets two new columns filled in for
#   the handed-back rows: E (Latest Target File) and F (Latest Handback File),
#   each a hyperlinked copy of the existing source-file / handoff-file links
# - Column G (Latest Handback DateTime) gets the real handback timestamp
#   instead of the "0001-01-01 00:00:00" placeholder

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Cells.Replace($oldStatus, $newStatus)
$wsZhCn.Cells.Replace($oldStatus, $newStatus)
$wsDeDe.Cells.Replace($oldStatus, $newStatus)

# ---------------- zh-cn sheet ----------------

# Row 2: 12866827-d98c-4a17-b73d-265334bbbe8d
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/02323729f9c8ab37803183d5b14f3ba48faadf38/e2e/12866827-d98c-4a17-b73d-265334bbbe8d.md", "", "", "12866827-d98c-4a17-b73d-265334bbbe8d.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7e323edafbf8a7561881dab4917633bb8fc2fbb4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/12866827-d98c-4a17-b73d-265334bbbe8d.79987c86519d0dbd026cdff0891d48d1f3e20ab6.zh-cn.xlf", "", "", "12866827-d98c-4a17-b73d-265334bbbe8d.79987c86519d0dbd026cdff0891d48d1f3e20ab6.zh-cn.xlf")
$wsZhCn.Range("G2").Value = "2016-01-26 12:26:01"

# Row 3: fb54192e-400e-44ac-b23e-5224e823a2da
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/02323729f9c8ab37803183d5b14f3ba48faadf38/e2e/fb54192e-400e-44ac-b23e-5224e823a2da.md", "", "", "fb54192e-400e-44ac-b23e-5224e823a2da.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7e323edafbf8a7561881dab4917633bb8fc2fbb4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/fb54192e-400e-44ac-b23e-5224e823a2da.e2299bd3a4ac64c525b1f5ed1fd64c9f7c101ddf.zh-cn.xlf", "", "", "fb54192e-400e-44ac-b23e-5224e823a2da.e2299bd3a4ac64c525b1f5ed1fd64c9f7c101ddf.zh-cn.xlf")
$wsZhCn.Range("G3").Value = "2016-01-26 12:26:01"

# ---------------- de-de sheet ----------------

# Row 2: 12866827-d98c-4a17-b73d-265334bbbe8d
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("E2"), "https://github.com/OpenLocalizationTest/oltest/blob/02323729f9c8ab37803183d5b14f3ba48faadf38/e2e/12866827-d98c-4a17-b73d-265334bbbe8d.md", "", "", "12866827-d98c-4a17-b73d-265334bbbe8d.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cb5a6b8848a1d4c35408af04a94c448bf0d6b028/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/12866827-d98c-4a17-b73d-265334bbbe8d.79987c86519d0dbd026cdff0891d48d1f3e20ab6.de-de.xlf", "", "", "12866827-d98c-4a17-b73d-265334bbbe8d.79987c86519d0dbd026cdff0891d48d1f3e20ab6.de-de.xlf")
$wsDeDe.Range("G2").Value = "2016-01-26 12:26:22"

# Row 3: fb54192e-400e-44ac-b23e-5224e823a2da
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("E3"), "https://github.com/OpenLocalizationTest/oltest/blob/02323729f9c8ab37803183d5b14f3ba48faadf38/e2e/fb54192e-400e-44ac-b23e-5224e823a2da.md", "", "", "fb54192e-400e-44ac-b23e-5224e823a2da.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cb5a6b8848a1d4c35408af04a94c448bf0d6b028/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/fb54192e-400e-44ac-b23e-5224e823a2da.e2299bd3a4ac64c525b1f5ed1fd64c9f7c101ddf.de-de.xlf", "", "", "fb54192e-400e-44ac-b23e-5224e823a2da.e2299bd3a4ac64c525b1f5ed1fd64c9f7c101ddf.de-de.xlf")
$wsDeDe.Range("G3").Value = "2016-01-26 12:26:22"
